$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 865.8076211556964
$ws.Range("B2").Value = 3.935489187071347
$ws.Range("C2").Value = 81.22553376900157
$ws.Range("D2").Value = 1.547982161899338
$ws.Range("E2").Value = 870.9627685546875
$ws.Range("F2").Value = 3.973341464996338
$ws.Range("G2").Value = 81.16103363037109
$ws.Range("H2").Value = 1.547192215919495
$ws.Range("A3").Value = 668.6751214703528
$ws.Range("B3").Value = 3.039432370319786
$ws.Range("C3").Value = 74.82121774012859
$ws.Range("D3").Value = 1.381594046610092
$ws.Range("E3").Value = 671.8682250976562
$ws.Range("F3").Value = 3.044100999832153
$ws.Range("G3").Value = 75.18991851806641
$ws.Range("H3").Value = 1.377697110176086
$ws.Range("A4").Value = 733.7113804798977
$ws.Range("B4").Value = 3.33505172945408
$ws.Range("C4").Value = 64.60832586860181
$ws.Range("D4").Value = 1.850540085054864
$ws.Range("E4").Value = 725.961669921875
$ws.Range("F4").Value = 3.302112817764282
$ws.Range("G4").Value = 64.74369049072266
$ws.Range("H4").Value = 1.840295076370239
$ws.Range("A5").Value = 745.9156777051711
$ws.Range("B5").Value = 3.390525807750778
$ws.Range("C5").Value = 76.14150837841709
$ws.Range("D5").Value = 1.40847087339762
$ws.Range("E5").Value = 755.6278686523438
$ws.Range("F5").Value = 3.416039705276489
$ws.Range("G5").Value = 76.54212951660156
$ws.Range("H5").Value = 1.418240070343018
$ws.Range("A6").Value = 403.7478290684457
$ws.Range("B6").Value = 1.835217404856571
$ws.Range("C6").Value = 55.12651421961891
$ws.Range("D6").Value = 1.22176114859273
$ws.Range("E6").Value = 408.8901062011719
$ws.Range("F6").Value = 1.86414647102356
$ws.Range("G6").Value = 54.97362518310547
$ws.Range("H6").Value = 1.226698160171509
$ws.Range("A7").Value = 707.2714450794919
$ws.Range("B7").Value = 3.214870204906781
$ws.Range("C7").Value = 69.81842864329923
$ws.Range("D7").Value = 1.614167055920128
$ws.Range("E7").Value = 699.8240356445312
$ws.Range("F7").Value = 3.169941186904907
$ws.Range("G7").Value = 69.64033508300781
$ws.Range("H7").Value = 1.603927612304688
$ws.Range("A8").Value = 667.1905554385579
$ws.Range("B8").Value = 3.032684342902536
$ws.Range("C8").Value = 56.70679239225603
$ws.Range("D8").Value = 1.87067564979227
$ws.Range("E8").Value = 661.5151977539062
$ws.Range("F8").Value = 3.029650688171387
$ws.Range("G8").Value = 56.74850463867188
$ws.Range("H8").Value = 1.872677326202393
$ws.Range("A9").Value = 747.4686885219467
$ws.Range("B9").Value = 3.39758494782703
$ws.Range("C9").Value = 79.15825930703772
$ws.Range("D9").Value = 1.441499830676714
$ws.Range("E9").Value = 768.3980712890625
$ws.Range("F9").Value = 3.455219507217407
$ws.Range("G9").Value = 79.20829772949219
$ws.Range("H9").Value = 1.461762070655823
$ws.Range("A10").Value = 801.8315002874823
$ws.Range("B10").Value = 3.644688637670374
$ws.Range("C10").Value = 65.51380157826003
$ws.Range("D10").Value = 2.086063669371808
$ws.Range("E10").Value = 819.1380004882812
$ws.Range("F10").Value = 3.778500556945801
$ws.Range("G10").Value = 65.61617279052734
$ws.Range("H10").Value = 2.133080244064331
$ws.Range("A11").Value = 522.8316059136538
$ws.Range("B11").Value = 2.376507299607517
$ws.Range("C11").Value = 64.98423315479118
$ws.Range("D11").Value = 1.210175877709597
$ws.Range("E11").Value = 521.57763671875
$ws.Range("F11").Value = 2.399527311325073
$ws.Range("G11").Value = 65.25044250488281
$ws.Range("H11").Value = 1.20848536491394
$ws.Range("A12").Value = 681.4826117151366
$ws.Range("B12").Value = 3.097648235068803
$ws.Range("C12").Value = 58.01969777608116
$ws.Range("D12").Value = 2.089637911197305
$ws.Range("E12").Value = 683.82861328125
$ws.Range("F12").Value = 3.141413450241089
$ws.Range("G12").Value = 58.33683395385742
$ws.Range("H12").Value = 2.087328672409058
$ws.Range("A13").Value = 912.6747679843742
$ws.Range("B13").Value = 4.148521672656246
$ws.Range("C13").Value = 74.96702938973311
$ws.Range("D13").Value = 1.813587795881309
$ws.Range("E13").Value = 915.4268188476562
$ws.Range("F13").Value = 4.16252613067627
$ws.Range("G13").Value = 74.97228240966797
$ws.Range("H13").Value = 1.812084317207336
$ws.Range("A14").Value = 849.4332832699217
$ws.Range("B14").Value = 3.861060378499644
$ws.Range("C14").Value = 64.03077297024208
$ws.Range("D14").Value = 2.205643060268716
$ws.Range("E14").Value = 859.3466186523438
$ws.Range("F14").Value = 3.909688711166382
$ws.Range("G14").Value = 64.09207153320312
$ws.Range("H14").Value = 2.212320327758789
$ws.Range("A15").Value = 764.9882203941365
$ws.Range("B15").Value = 3.477219183609712
$ws.Range("C15").Value = 67.58879894320933
$ws.Range("D15").Value = 1.847690605435842
$ws.Range("E15").Value = 773.1941528320312
$ws.Range("F15").Value = 3.512359380722046
$ws.Range("G15").Value = 67.32568359375
$ws.Range("H15").Value = 1.879199504852295
$ws.Range("A16").Value = 646.8230595060777
$ws.Range("B16").Value = 2.940104815936717
$ws.Range("C16").Value = 54.72920047364931
$ws.Range("D16").Value = 1.738365167247628
$ws.Range("E16").Value = 661.9425659179688
$ws.Range("F16").Value = 2.991424083709717
$ws.Range("G16").Value = 54.3863639831543
$ws.Range("H16").Value = 1.742787003517151
$ws.Range("A17").Value = 665.6079907264383
$ws.Range("B17").Value = 3.025490866938356
$ws.Range("C17").Value = 73.26153166599943
$ws.Range("D17").Value = 1.222264593286663
$ws.Range("E17").Value = 695.0214233398438
$ws.Range("F17").Value = 3.153612375259399
$ws.Range("G17").Value = 73.56150054931641
$ws.Range("H17").Value = 1.238518238067627
$ws.Range("A18").Value = 779.9062136581094
$ws.Range("B18").Value = 3.545028243900497
$ws.Range("C18").Value = 71.98595204117288
$ws.Range("D18").Value = 1.507107925782621
$ws.Range("E18").Value = 820.4475708007812
$ws.Range("F18").Value = 3.666774272918701
$ws.Range("G18").Value = 72.11370849609375
$ws.Range("H18").Value = 1.561562776565552
$ws.Range("A19").Value = 564.8398388976129
$ws.Range("B19").Value = 2.567453813170968
$ws.Range("C19").Value = 69.18527132628077
$ws.Range("D19").Value = 1.309963659524748
$ws.Range("E19").Value = 575.35400390625
$ws.Range("F19").Value = 2.664212703704834
$ws.Range("G19").Value = 69.37921905517578
$ws.Range("H19").Value = 1.343511462211609
$ws.Range("A20").Value = 869.5931556945404
$ws.Range("B20").Value = 3.952696162247911
$ws.Range("C20").Value = 78.77910442401628
$ws.Range("D20").Value = 1.703559075869232
$ws.Range("E20").Value = 873.0420532226562
$ws.Range("F20").Value = 3.946182489395142
$ws.Range("G20").Value = 78.64844512939453
$ws.Range("H20").Value = 1.694158673286438
$ws.Range("A21").Value = 592.9610695625362
$ws.Range("B21").Value = 2.695277588920619
$ws.Range("C21").Value = 70.99553502428301
$ws.Range("D21").Value = 1.335685312556614
$ws.Range("E21").Value = 607.6273193359375
$ws.Range("F21").Value = 2.806180953979492
$ws.Range("G21").Value = 71.07904052734375
$ws.Range("H21").Value = 1.367692947387695
